$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 15582.8
$ws.Range("I33").Value = 15582.8
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 15582.8
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -15353.8
$ws.Range("N33").ClearContents()
$ws.Range("H58").Value = 7629.857
$ws.Range("J58").Value = 1500
$ws.Range("L58").Value = 4500
$ws.Range("N58").Value = -4800
$ws.Range("H74").Value = 9917.166999999999
$ws.Range("I74").Value = 9751.5
$ws.Range("K74").Value = 9751.5
$ws.Range("M74").Value = -8815.5
$ws.Range("H76").Value = 20009388
$ws.Range("I76").Value = 33337998
$ws.Range("K76").Value = 33337998
$ws.Range("M76").Value = -33337683
$ws.Range("H77").Value = 9917.166999999999
$ws.Range("I77").Value = 9751.5
$ws.Range("K77").Value = 48757.5
$ws.Range("M77").Value = -44077.5
$ws.Range("H79").Value = 20009388
$ws.Range("I79").Value = 33337998
$ws.Range("K79").Value = 33337998
$ws.Range("M79").Value = -33336906
$ws.Range("H116").Value = 3849
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H137").Value = 9620302
$ws.Range("I137").Value = 15626134
$ws.Range("J137").Value = 10971.2
$ws.Range("K137").Value = 46878402
$ws.Range("L137").Value = 32913.60000000001
$ws.Range("M137").Value = -46875852
$ws.Range("N137").Value = -38013.60000000001
$ws.Range("H141").Value = 984.6667
$ws.Range("I141").Value = 984.6667
$ws.Range("K141").Value = 2954.0001
$ws.Range("M141").Value = 2225.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3361
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -885
$ws.Range("H4").Value = 490
$ws.Range("I4").Value = 362.5
$ws.Range("K4").Value = 362.5
$ws.Range("M4").Value = -246.5
$ws.Range("H32").Value = 20679.424
$ws.Range("I32").Value = 20496.666
$ws.Range("K32").Value = 20496.666
$ws.Range("M32").Value = -20209.666
$ws.Range("H45").Value = 8649.5
$ws.Range("I45").Value = 10263
$ws.Range("K45").Value = 10263
$ws.Range("M45").Value = -9886
$ws.Range("H132").Value = 3174.6155
$ws.Range("I132").Value = 2097.3333
$ws.Range("J132").Value = 7699.2
$ws.Range("K132").Value = 6291.999899999999
$ws.Range("L132").Value = 23097.6
$ws.Range("M132").Value = -3761.999899999999
$ws.Range("N132").Value = -28157.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 795.3684
$ws.Range("I94").Value = 758.3570999999999
$ws.Range("J94").Value = 899
$ws.Range("K94").Value = 758.3570999999999
$ws.Range("L94").Value = 899
$ws.Range("M94").Value = -307.3570999999999
$ws.Range("N94").Value = -1801
$ws.Range("H134").Value = 4996.1665
$ws.Range("I134").Value = 2413.6
$ws.Range("K134").Value = 7240.799999999999
$ws.Range("M134").Value = -4705.799999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32263356
$ws.Range("I31").Value = 90910160
$ws.Range("J31").Value = 7612.85
$ws.Range("K31").Value = 90910160
$ws.Range("L31").Value = 7612.85
$ws.Range("M31").Value = -90909865
$ws.Range("N31").Value = -8202.85
$ws.Range("H34").Value = 32263356
$ws.Range("I34").Value = 90910160
$ws.Range("J34").Value = 7612.85
$ws.Range("K34").Value = 90910160
$ws.Range("L34").Value = 7612.85
$ws.Range("M34").Value = -90909958
$ws.Range("N34").Value = -8016.85
$ws.Range("H58").Value = 4659.625
$ws.Range("I58").Value = 3197.55
$ws.Range("J58").Value = 7096.4165
$ws.Range("K58").Value = 3197.55
$ws.Range("L58").Value = 7096.4165
$ws.Range("M58").Value = -2994.55
$ws.Range("N58").Value = -7502.4165
$ws.Range("H99").Value = 1951
$ws.Range("I99").Value = 1001
$ws.Range("J99").Value = 2901
$ws.Range("K99").Value = 1001
$ws.Range("L99").Value = 2901
$ws.Range("M99").Value = 497
$ws.Range("N99").Value = -5897
$ws.Range("H122").Value = 63723.375
$ws.Range("I122").Value = 100969.7
$ws.Range("J122").Value = 1646.1666
$ws.Range("K122").Value = 302909.1
$ws.Range("L122").Value = 4938.4998
$ws.Range("M122").Value = -300459.1
$ws.Range("N122").Value = -9838.4998
$ws.Range("H126").Value = 1951
$ws.Range("I126").Value = 1001
$ws.Range("J126").Value = 2901
$ws.Range("K126").Value = 3003
$ws.Range("L126").Value = 8703
$ws.Range("M126").Value = -533
$ws.Range("N126").Value = -13643
$ws.Range("H136").Value = 4659.625
$ws.Range("I136").Value = 3197.55
$ws.Range("J136").Value = 7096.4165
$ws.Range("K136").Value = 9592.650000000001
$ws.Range("L136").Value = 21289.2495
$ws.Range("M136").Value = -7042.650000000001
$ws.Range("N136").Value = -26389.2495
$ws.Range("H138").Value = 120000
$ws.Range("J138").Value = 120000
$ws.Range("L138").Value = 120000
$ws.Range("N138").Value = -130280
$ws.Range("H139").Value = 90139.5
$ws.Range("J139").Value = 90139.5
$ws.Range("L139").Value = 90139.5
$ws.Range("N139").Value = -100419.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2002.25
$ws.Range("J25").Value = 2503.3333
$ws.Range("L25").Value = 7509.999899999999
$ws.Range("N25").Value = -7847.999899999999
$ws.Range("H30").Value = 2002.25
$ws.Range("J30").Value = 2503.3333
$ws.Range("L30").Value = 7509.999899999999
$ws.Range("N30").Value = -7713.999899999999
$ws.Range("H136").Value = 2588.6667
$ws.Range("I136").Value = 366.5
$ws.Range("K136").Value = 1099.5
$ws.Range("M136").Value = 4000.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 671.619
$ws.Range("I2").Value = 949.6429000000001
$ws.Range("J2").Value = 115.57143
$ws.Range("K2").Value = 949.6429000000001
$ws.Range("L2").Value = 115.57143
$ws.Range("M2").Value = -836.6429000000001
$ws.Range("N2").Value = -341.57143
$ws.Range("H70").Value = 29667.846
$ws.Range("I70").Value = 25961.75
$ws.Range("J70").Value = 35597.6
$ws.Range("K70").Value = 25961.75
$ws.Range("L70").Value = 35597.6
$ws.Range("M70").Value = -25691.75
$ws.Range("N70").Value = -36137.6
$ws.Range("H73").Value = 29667.846
$ws.Range("I73").Value = 25961.75
$ws.Range("J73").Value = 35597.6
$ws.Range("K73").Value = 25961.75
$ws.Range("L73").Value = 35597.6
$ws.Range("M73").Value = -25025.75
$ws.Range("N73").Value = -37469.6
$ws.Range("H126").Value = 3249.606
$ws.Range("J126").Value = 4102.533
$ws.Range("L126").Value = 12307.599
$ws.Range("N126").Value = -17247.599
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 44499.332
$ws.Range("J136").Value = 44499.332
$ws.Range("L136").Value = 133497.996
$ws.Range("N136").Value = -138597.996
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4386.9165
$ws.Range("J46").Value = 4720.3257
$ws.Range("L46").Value = 4720.3257
$ws.Range("N46").Value = -5096.3257
$ws.Range("H132").Value = 4554.615
$ws.Range("I132").Value = 3421.2222
$ws.Range("J132").Value = 7104.75
$ws.Range("K132").Value = 10263.6666
$ws.Range("L132").Value = 21314.25
$ws.Range("M132").Value = -7733.6666
$ws.Range("N132").Value = -26374.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7171.4287
$ws.Range("I81").Value = 1968.5
$ws.Range("J81").Value = 10373.23
$ws.Range("K81").Value = 3937
$ws.Range("L81").Value = 20746.46
$ws.Range("M81").Value = -2876
$ws.Range("N81").Value = -22868.46
$ws.Range("H84").Value = 7171.4287
$ws.Range("I84").Value = 1968.5
$ws.Range("J84").Value = 10373.23
$ws.Range("K84").Value = 19685
$ws.Range("L84").Value = 103732.3
$ws.Range("M84").Value = -14381
$ws.Range("N84").Value = -114340.3
$ws.Range("H123").Value = 49999.5
$ws.Range("J123").Value = 49999.5
$ws.Range("L123").Value = 49999.5
$ws.Range("N123").Value = -59799.5
